# Gamelab_2_Planning.xlsx edit script
# Commit: "A* en shaders update (not finished yet)"
#
# Summary of the change:
#  - Sheet "Blad1": a new "Shaders" section (3 rows) is inserted right after the
#    "Props" section (after row 83), pushing the "Particles" section down.
#    The "Bleeding" particle entry is removed from the Particles section.
#  - Sheet "Blad2": a couple of cells that (by coincidence of shared-string
#    reuse) used to show "PAR_BLEED" now show "Skin shader" / stay "PAR_BREATH".
#  - Active sheet / selection moves from Blad2 to Blad1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Insert 3 new rows for the "Shaders" section, right below "Props" (row 83)
# ---------------------------------------------------------------------------
$ws1.Rows("84:86").Insert()

# New row 84: section header "Shaders" (style copied from the "Props" header, A81)
$ws1.Range("A81").Copy()
$ws1.Range("A84").PasteSpecial(-4122)
$ws1.Range("A84").Value = "Shaders"

# B84 picked up the (accidental) "Opmerkingen" header fill in the original edit
$ws1.Range("H1").Copy()
$ws1.Range("B84").PasteSpecial(-4122)

$ws1.Range("H5").Copy()
$ws1.Range("H84").PasteSpecial(-4122)

$ws1.Range("J81").Copy()
$ws1.Range("J84").PasteSpecial(-4122)
$ws1.Range("I84").ClearContents()
$ws1.Range("I84").ClearFormats()

# New row 85: "Skin shader" / remark "Makes blood appear on models" / owner Alieke
$ws1.Range("J81").Copy()
$ws1.Range("A85").PasteSpecial(-4122)

$ws1.Range("B82").Copy()
$ws1.Range("B85").PasteSpecial(-4122)
$ws1.Range("B85").Value = "Skin shader"

$ws1.Range("H12").Copy()
$ws1.Range("H85").PasteSpecial(-4122)
$ws1.Range("H85").Value = "Makes blood appear on models"

$ws1.Range("I82").Copy()
$ws1.Range("I85").PasteSpecial(-4122)
$ws1.Range("I85").Value = "Alieke"

# New row 86: "Mos shader"
$ws1.Range("J81").Copy()
$ws1.Range("A86").PasteSpecial(-4122)

$ws1.Range("B82").Copy()
$ws1.Range("B86").PasteSpecial(-4122)
$ws1.Range("B86").Value = "Mos shader"

$ws1.Range("H5").Copy()
$ws1.Range("H86").PasteSpecial(-4122)

$ws1.Range("J81").Copy()
$ws1.Range("J86").PasteSpecial(-4122)
$ws1.Range("I86").ClearContents()
$ws1.Range("I86").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Remove the "Bleeding" particle entry (was row 88, now row 91 after insert)
#    by shifting the rest of the Particles rows up by one.
# ---------------------------------------------------------------------------
$ws1.Rows("91").Delete()

# ---------------------------------------------------------------------------
# 3. Sheet "Blad2": text content tweaks
# ---------------------------------------------------------------------------
$ws2.Range("G3").Value = "Skin shader"

# ---------------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping (cosmetic, matches author's last view)
# ---------------------------------------------------------------------------
$ws1.Select()
$excel.ActiveWindow.ScrollRow = 61
$ws1.Range("B84").Select()
